$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 311, shifting existing rows 311-426 down to 312-427
$ws.Rows("311:311").Insert()

# Fill in the new row 311 with its data (same as the row that was below it, except for the
# columns that actually changed value)
$ws.Range("A311").Value = 9
$ws.Range("B311").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C311").Value = "Metropolitana"
$ws.Range("D311").Value = 44795
$ws.Range("E311").Value = 13
$ws.Range("F311").Value = 100112044
$ws.Range("G311").Value = "Perejil"
$ws.Range("H311").Value = "Sin especificar"
$ws.Range("I311").Value = "Primera"
$ws.Range("J311").Value = 52
$ws.Range("K311").Value = 17000
$ws.Range("L311").Value = 18000
$ws.Range("M311").Value = 17500
$ws.Range("N311").Value = "$/docena de atados"
$ws.Range("O311").Value = "Región Metropolitana"
$ws.Range("P311").Value = 5833
$ws.Range("Q311").Value = 3
$ws.Range("R311").Value = "Hortaliza"

$ws.Range("D311").NumberFormat = "YYYY-MM-DD HH:MM:SS"
